$wb = $excel.ActiveWorkbook

# "Generate Report for Handoff" — refresh the handoff timestamps for the
# ede29265-f53d-414c-9b14-52b38b8de1ae row across the Overview, zh-cn and
# de-de sheets.

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G7").Value = "2016-09-04 20:47:37"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H7").Value = "2016-09-04 20:47:32"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H7").Value = "2016-09-04 20:47:37"
